$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert the new "NEW GRAPH" list item (+ trailing blank spacer paragraph)
#    right after the "Notification for when a new Crypto is created -> ..."
#    bullet, i.e. immediately before the first of the run of blank
#    ListParagraph paragraphs that precede the "Monetization:" title.
# ---------------------------------------------------------------------------
$anchor = $d.Content
$anchor.Find.Execute("Notification for when a new Crypto is created") | Out-Null
$anchor.Expand(4) | Out-Null                 # wdParagraph -> whole bullet paragraph
$blank = $anchor.Next(4, 1)                  # the first blank paragraph right after it
$insertionPoint = $d.Range($blank.End, $blank.End)

$newGraphXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr>
<w:pStyle w:val="ListParagraph"/>
<w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>
<w:rPr><w:b/><w:bCs/></w:rPr>
</w:pPr>
<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">NEW GRAPH: </w:t></w:r>
<w:r><w:t>Relation in DB that tracks how any views/votes/comments a stock gets per day, regardless of post_id, CREATE TABLE (symbol, date, votes, views, comments), PRIMARY KEY(symbol, date)</w:t></w:r>
</w:p>
<w:p>
<w:pPr>
<w:pStyle w:val="ListParagraph"/>
<w:rPr><w:b/><w:bCs/></w:rPr>
</w:pPr>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$insertionPoint.InsertXML($newGraphXml)

# ---------------------------------------------------------------------------
# 2) Collapse the two "Monetization" / ":" runs (and drop the stray
#    lastRenderedPageBreak) into a single run reading "Monetization:".
# ---------------------------------------------------------------------------
$title = $d.Content
$title.Find.Execute("Monetization") | Out-Null
$title.Expand(4) | Out-Null                  # wdParagraph -> whole title paragraph

$titleXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr><w:pStyle w:val="Title"/></w:pPr>
<w:r><w:t>Monetization:</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$title.InsertXML($titleXml)

Write-Output "Applied IDEAS.docx edit: added NEW GRAPH bullet + merged Monetization title runs."
